$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data gained a new weekly record. It belongs logically right
# before the existing row 99, so insert a blank row there (this pushes the
# former rows 99-130 down to 100-131, preserving all of their data/styles).
$ws.Rows.Item(99).Insert()

# Fill the newly inserted row 99 with the new "Apio" market record.
$ws.Cells.Item(99, 1).Value2  = 4
$ws.Cells.Item(99, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(99, 3).Value2  = "Los Lagos"
$ws.Cells.Item(99, 4).Value2  = 44463
$ws.Cells.Item(99, 5).Value2  = 10
$ws.Cells.Item(99, 6).Value2  = 100112017
$ws.Cells.Item(99, 7).Value2  = "Apio"
$ws.Cells.Item(99, 8).Value2  = "Americana (o)"
$ws.Cells.Item(99, 9).Value2  = "Primera"
$ws.Cells.Item(99, 10).Value2 = 40
$ws.Cells.Item(99, 11).Value2 = 14000
$ws.Cells.Item(99, 12).Value2 = 14000
$ws.Cells.Item(99, 13).Value2 = 14000
$ws.Cells.Item(99, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(99, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(99, 16).Value2 = 2333
$ws.Cells.Item(99, 17).Value2 = 6
$ws.Cells.Item(99, 18).Value2 = "Hortaliza"
